# Move heatshield analysis into production:
#  - change the loss/diss exponents in K8:M8
#  - add a new "abl/area" column (J) and a second loss/diss block (O:Q)
#    expressed per unit mass instead of per diameter
#  - update the view (scroll position / selection)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- updated exponents used by the existing K:M formulas -------------
$ws.Range("K8").Value = 2.25
$ws.Range("L8").Value = 2.25
$ws.Range("M8").Value = 2

# --- new header row cells --------------------------------------------
$ws.Range("J1").Value = "abl/area"
$ws.Range("O1").Value = "loss1000"
$ws.Range("P1").Value = "loss3000"
$ws.Range("Q1").Value = "diss500"

# --- new "per area" block, rows 2-5 -----------------------------------
$ws.Range("O2").Formula = '=F2/$E2'
$ws.Range("P2").Formula = '=G2/$E2'
$ws.Range("Q2").Formula = '=H2*$E2'

$ws.Range("O3:O5").Formula = '=F3/$E3'
$ws.Range("P3:P5").Formula = '=G3/$E3'
$ws.Range("Q3:Q5").Formula = '=H3*$E3'

# --- view: scroll/selection match the published state -----------------
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("Q4").Select()
